# Insert a new data row above row 175, shifting existing rows 175-251 down to 176-252.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(175).Insert()

$ws.Cells.Item(175, 1).Value2  = 9
$ws.Cells.Item(175, 2).Value   = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(175, 3).Value   = "Metropolitana"
$ws.Cells.Item(175, 4).Value2  = 44627
$ws.Cells.Item(175, 5).Value2  = 13
$ws.Cells.Item(175, 6).Value2  = 100112001
$ws.Cells.Item(175, 7).Value   = "Berenjena"
$ws.Cells.Item(175, 8).Value   = "Sin especificar"
$ws.Cells.Item(175, 9).Value   = "Primera"
$ws.Cells.Item(175, 10).Value2 = 79
$ws.Cells.Item(175, 11).Value2 = 7000
$ws.Cells.Item(175, 12).Value2 = 8000
$ws.Cells.Item(175, 13).Value2 = 7494
$ws.Cells.Item(175, 14).Value  = "`$/caja 60 unidades"
$ws.Cells.Item(175, 15).Value  = "Región de Arica y Parinacota"
$ws.Cells.Item(175, 16).Value2 = 125
$ws.Cells.Item(175, 17).Value2 = 60
$ws.Cells.Item(175, 18).Value  = "Hortaliza"
